$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "DBLP Conference" header column (E) + fill existing rows where a DBLP
# journal abbreviation already exists, plus "-" placeholders in column D
# (DBLP Journal) for rows that don't have one.

$ws.Range("E1").Value = "DBLP Conference"
$ws.Range("E2").Value = "sigmod"
$ws.Range("E3").Value = "kdd"
$ws.Range("E4").Value = "sigir"
$ws.Range("E5").Value = "cvpr"

$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = "nips"

$ws.Range("E7").Value = "sigcomm"

$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "ccs"

$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = "icse"

$ws.Range("E10").Value = "isca"

$ws.Range("E11").Value = "chi"
$ws.Range("E12").Value = "podc"
$ws.Range("E13").Value = "siggraph"
$ws.Range("E14").Value = "recomb"
$ws.Range("E15").Value = "mm"

$ws.Columns.Item(4).ColumnWidth = 11.3984375
$ws.Columns.Item(5).ColumnWidth = 9.46484375

$ws.Range("D10").Select() | Out-Null
